# Update the "Förändrad" (changed) date column (C) for all data rows
# (rows 2 through 108) from serial date 45175 (2023-09-06) to
# serial date 45183 (2023-09-14), keeping the existing cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 108; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
